$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.92"
$ws.Range("D3").Value = "'22.03"
$ws.Range("D4").Value = "'5.440"
$ws.Range("D5").Value = "'0.05777"
$ws.Range("D8").Value = "'0.8163"
$ws.Range("D9").Value = "'1.045"
$ws.Range("E9").Value = '8FTXTokenFTTBestin24h'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1426"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").Value = "'0.07298"
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").Value = "'0.03110"
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03121"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = "'4.153"
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.09384"
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001591"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04808"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005847"
$ws.Range("E18").Value = '17OneONE'
$ws.Range("D19").Value = "'0.006276"
$ws.Range("D20").Value = "'0.004126"
$ws.Range("D21").Value = "'0.0009934"
$ws.Range("D23").Value = "'3.744"
$ws.Range("D24").Value = "'2.185"
$ws.Range("D26").Value = "'0.1329"
$ws.Range("D27").Value = "'0.0003997"
$ws.Range("D40").Value = "'0.03870"
$ws.Range("D41").Value = "'0.006694"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("D43").Value = "'0.002638"
$ws.Range("D44").Value = "'0.006570"
$ws.Range("D45").Value = "'0.00005614"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.3898"
$ws.Range("D49").Value = "'0.00002099"
